$d = $word.ActiveDocument

# The paragraph that currently only contains the "_GoBack" bookmark needs
# three new runs added around the bookmark: "TTTclass:" + "importance"
# before it, and "TTT" after it. Two new, empty paragraphs are then added
# right after that paragraph (before the final section properties).

$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range

# Insert the "TTT" run *after* the bookmark first (so the bookmark's
# position is not disturbed by the text inserted ahead of it later).
$afterRange = $bmRange.Duplicate
$afterRange.Collapse(0)   # wdCollapseEnd
$afterRange.InsertAfter("TTT")

# Insert the "TTTclass:" and "importance" runs *before* the bookmark.
$bm = $d.Bookmarks.Item("_GoBack")
$beforeRange = $bm.Range.Duplicate
$beforeRange.Collapse(1)  # wdCollapseStart
$beforeRange.InsertBefore("importance")
$beforeRange.InsertBefore("TTTclass:")

# The bookmark's paragraph is now the last paragraph in the body (the
# section-properties sentinel paragraph is not counted in .Paragraphs).
$bmParaIndex = $d.Paragraphs.Count
$bmPara = $d.Paragraphs.Item($bmParaIndex)

# Append two new paragraphs right after it.
$tail = $bmPara.Range.Duplicate
$tail.Collapse(0)  # wdCollapseEnd
$tail.InsertParagraphAfter()
$tail.InsertParagraphAfter()

# Freshly-created paragraph marks get serialised with a stray empty run;
# rewrite each new paragraph's OOXML to a clean, empty <w:p/> so it
# matches the shape of the other blank paragraphs already in the
# document.
$emptyParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$d.Paragraphs.Item($bmParaIndex + 1).Range.InsertXML($emptyParaXml)
$d.Paragraphs.Item($bmParaIndex + 2).Range.InsertXML($emptyParaXml)
